$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'ECs'
$ws.Range("B2").Value = 'Ntrk3'
$ws.Range("C2").Value = 'Ptprf'
$ws.Range("D2").Value = 'ECs'
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.002145666666666667
$ws.Range("H2").Value = 0.006437
$ws.Range("I2").Value = 0.00807320947388686
$ws.Range("J2").Value = 0.008843624333499573
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1255626666666667
$ws.Range("N2").Value = 0.376688
$ws.Range("O2").Value = 0.02744849445093922
$ws.Range("P2").Value = 0.02907345870642374
$ws.Range("Q2").Value = 0.0002694156284444445
$ws.Range("R2").Value = 0.002424740656
$ws.Range("S2").Value = 0.0002215974454452534
$ws.Range("T2").Value = 0.000257114746875124

$ws.Range("A3").Value = 'ECs'
$ws.Range("B3").Value = 'Ntrk3'
$ws.Range("C3").Value = 'Ptprf'
$ws.Range("D3").Value = 'FAPs'
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.002145666666666667
$ws.Range("H3").Value = 0.006437
$ws.Range("I3").Value = 0.00807320947388686
$ws.Range("J3").Value = 0.008843624333499573
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.393572666666666
$ws.Range("N3").Value = 10.180718
$ws.Range("O3").Value = 0.7418483772500767
$ws.Range("P3").Value = 0.7857661629113346
$ws.Range("Q3").Value = 0.007281475751777777
$ws.Range("R3").Value = 0.065533281766
$ws.Range("S3").Value = 0.005989097347402913
$ws.Range("T3").Value = 0.006949020758763269

$ws.Range("A4").Value = 'ECs'
$ws.Range("B4").Value = 'Ntrk3'
$ws.Range("C4").Value = 'Ptprf'
$ws.Range("D4").Value = 'Inflammatory-Mac'
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.002145666666666667
$ws.Range("H4").Value = 0.006437
$ws.Range("I4").Value = 0.00807320947388686
$ws.Range("J4").Value = 0.008843624333499573
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2723486666666667
$ws.Range("N4").Value = 0.8170460000000001
$ws.Range("O4").Value = 0.05953649332381727
$ws.Range("P4").Value = 0.06306108275880487
$ws.Range("Q4").Value = 0.0005843694557777779
$ws.Range("R4").Value = 0.005259325102
$ws.Range("S4").Value = 0.0004806505819438434
$ws.Range("T4").Value = 0.0005576885259825972

$ws.Range("A5").Value = 'ECs'
$ws.Range("B5").Value = 'Ntrk3'
$ws.Range("C5").Value = 'Ptprf'
$ws.Range("D5").Value = 'MuSCs'
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.002145666666666667
$ws.Range("H5").Value = 0.006437
$ws.Range("I5").Value = 0.00807320947388686
$ws.Range("J5").Value = 0.008843624333499573
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7670265000000001
$ws.Range("N5").Value = 1.534053
$ws.Range("O5").Value = 0.167675019875653
$ws.Range("P5").Value = 0.1184009752075072
$ws.Range("Q5").Value = 0.0016457831935
$ws.Range("R5").Value = 0.009874699161
$ws.Range("S5").Value = 0.001353675558994289
$ws.Range("T5").Value = 0.00104709374545519

$ws.Range("A6").Value = 'ECs'
$ws.Range("B6").Value = 'Ntrk3'
$ws.Range("C6").Value = 'Ptprf'
$ws.Range("D6").Value = 'Resolving-Mac'
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.002145666666666667
$ws.Range("H6").Value = 0.006437
$ws.Range("I6").Value = 0.00807320947388686
$ws.Range("J6").Value = 0.008843624333499573
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01597233333333333
$ws.Range("N6").Value = 0.047917
$ws.Range("O6").Value = 0.0034916150995138
$ws.Range("P6").Value = 0.003698320415929645
$ws.Range("Q6").Value = 0.00003427130322222223
$ws.Range("R6").Value = 0.000308441729
$ws.Range("S6").Value = 0.00002818854010056122
$ws.Range("T6").Value = 0.00003270655642339368

$ws.Range("A7").Value = 'FAPs'
$ws.Range("B7").Value = 'Ntrk3'
$ws.Range("C7").Value = 'Ptprf'
$ws.Range("D7").Value = 'ECs'
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.03558433333333334
$ws.Range("H7").Value = 0.106753
$ws.Range("I7").Value = 0.1338883534202026
$ws.Range("J7").Value = 0.1466651279282398
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.1255626666666667
$ws.Range("N7").Value = 0.376688
$ws.Range("O7").Value = 0.02744849445093922
$ws.Range("P7").Value = 0.02907345870642374
$ws.Range("Q7").Value = 0.00446806378488889
$ws.Range("R7").Value = 0.040212574064
$ws.Range("S7").Value = 0.003675033725899819
$ws.Range("T7").Value = 0.004264062540494036

$ws.Range("A8").Value = 'FAPs'
$ws.Range("B8").Value = 'Ntrk3'
$ws.Range("C8").Value = 'Ptprf'
$ws.Range("D8").Value = 'FAPs'
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.03558433333333334
$ws.Range("H8").Value = 0.106753
$ws.Range("I8").Value = 0.1338883534202026
$ws.Range("J8").Value = 0.1466651279282398
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.393572666666666
$ws.Range("N8").Value = 10.180718
$ws.Range("O8").Value = 0.7418483772500767
$ws.Range("P8").Value = 0.7857661629113346
$ws.Range("Q8").Value = 0.1207580209615556
$ws.Range("R8").Value = 1.086822188654
$ws.Range("S8").Value = 0.09932485771746204
$ws.Range("T8").Value = 0.115244494805073

$ws.Range("A9").Value = 'FAPs'
$ws.Range("B9").Value = 'Ntrk3'
$ws.Range("C9").Value = 'Ptprf'
$ws.Range("D9").Value = 'Inflammatory-Mac'
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.03558433333333334
$ws.Range("H9").Value = 0.106753
$ws.Range("I9").Value = 0.1338883534202026
$ws.Range("J9").Value = 0.1466651279282398
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2723486666666667
$ws.Range("N9").Value = 0.8170460000000001
$ws.Range("O9").Value = 0.05953649332381727
$ws.Range("P9").Value = 0.06306108275880487
$ws.Range("Q9").Value = 0.009691345737555557
$ws.Range("R9").Value = 0.087222111638
$ws.Range("S9").Value = 0.00797124305953878
$ws.Range("T9").Value = 0.009248861770113435

$ws.Range("A10").Value = 'FAPs'
$ws.Range("B10").Value = 'Ntrk3'
$ws.Range("C10").Value = 'Ptprf'
$ws.Range("D10").Value = 'MuSCs'
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.03558433333333334
$ws.Range("H10").Value = 0.106753
$ws.Range("I10").Value = 0.1338883534202026
$ws.Range("J10").Value = 0.1466651279282398
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.7670265000000001
$ws.Range("N10").Value = 1.534053
$ws.Range("O10").Value = 0.167675019875653
$ws.Range("P10").Value = 0.1184009752075072
$ws.Range("Q10").Value = 0.02729412665150001
$ws.Range("R10").Value = 0.163764759909
$ws.Range("S10").Value = 0.02244973232085092
$ws.Range("T10").Value = 0.01736529417563739

$ws.Range("A11").Value = 'FAPs'
$ws.Range("B11").Value = 'Ntrk3'
$ws.Range("C11").Value = 'Ptprf'
$ws.Range("D11").Value = 'Resolving-Mac'
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.03558433333333334
$ws.Range("H11").Value = 0.106753
$ws.Range("I11").Value = 0.1338883534202026
$ws.Range("J11").Value = 0.1466651279282398
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.01597233333333333
$ws.Range("N11").Value = 0.047917
$ws.Range("O11").Value = 0.0034916150995138
$ws.Range("P11").Value = 0.003698320415929645
$ws.Range("Q11").Value = 0.0005683648334444446
$ws.Range("R11").Value = 0.005115283501
$ws.Range("S11").Value = 0.0004674865964510195
$ws.Range("T11").Value = 0.0005424146369219426

$ws.Range("A12").Value = 'Inflammatory-Mac'
$ws.Range("B12").Value = 'Ntrk3'
$ws.Range("C12").Value = 'Ptprf'
$ws.Range("D12").Value = 'ECs'
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.04573
$ws.Range("H12").Value = 0.13719
$ws.Range("I12").Value = 0.172062079807758
$ws.Range("J12").Value = 0.188481718551003
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.1255626666666667
$ws.Range("N12").Value = 0.376688
$ws.Range("O12").Value = 0.02744849445093922
$ws.Range("P12").Value = 0.02907345870642374
$ws.Range("Q12").Value = 0.005741980746666668
$ws.Range("R12").Value = 0.05167782672000001
$ws.Range("S12").Value = 0.004722845042820306
$ws.Range("T12").Value = 0.005479815461208367

$ws.Range("A13").Value = 'Inflammatory-Mac'
$ws.Range("B13").Value = 'Ntrk3'
$ws.Range("C13").Value = 'Ptprf'
$ws.Range("D13").Value = 'FAPs'
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.04573
$ws.Range("H13").Value = 0.13719
$ws.Range("I13").Value = 0.172062079807758
$ws.Range("J13").Value = 0.188481718551003
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.393572666666666
$ws.Range("N13").Value = 10.180718
$ws.Range("O13").Value = 0.7418483772500767
$ws.Range("P13").Value = 0.7857661629113346
$ws.Range("Q13").Value = 0.1551880780466666
$ws.Range("R13").Value = 1.39669270242
$ws.Range("S13").Value = 0.1276439746916585
$ws.Range("T13").Value = 0.1481025567647558

$ws.Range("A14").Value = 'Inflammatory-Mac'
$ws.Range("B14").Value = 'Ntrk3'
$ws.Range("C14").Value = 'Ptprf'
$ws.Range("D14").Value = 'Inflammatory-Mac'
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.04573
$ws.Range("H14").Value = 0.13719
$ws.Range("I14").Value = 0.172062079807758
$ws.Range("J14").Value = 0.188481718551003
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.2723486666666667
$ws.Range("N14").Value = 0.8170460000000001
$ws.Range("O14").Value = 0.05953649332381727
$ws.Range("P14").Value = 0.06306108275880487
$ws.Range("Q14").Value = 0.01245450452666667
$ws.Range("R14").Value = 0.11209054074
$ws.Range("S14").Value = 0.0102439728657567
$ws.Range("T14").Value = 0.01188586125206657

$ws.Range("A15").Value = 'Inflammatory-Mac'
$ws.Range("B15").Value = 'Ntrk3'
$ws.Range("C15").Value = 'Ptprf'
$ws.Range("D15").Value = 'MuSCs'
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.04573
$ws.Range("H15").Value = 0.13719
$ws.Range("I15").Value = 0.172062079807758
$ws.Range("J15").Value = 0.188481718551003
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.7670265000000001
$ws.Range("N15").Value = 1.534053
$ws.Range("O15").Value = 0.167675019875653
$ws.Range("P15").Value = 0.1184009752075072
$ws.Range("Q15").Value = 0.035076121845
$ws.Range("R15").Value = 0.21045673107
$ws.Range("S15").Value = 0.02885051265161202
$ws.Range("T15").Value = 0.02231641928522565

$ws.Range("A16").Value = 'Inflammatory-Mac'
$ws.Range("B16").Value = 'Ntrk3'
$ws.Range("C16").Value = 'Ptprf'
$ws.Range("D16").Value = 'Resolving-Mac'
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.04573
$ws.Range("H16").Value = 0.13719
$ws.Range("I16").Value = 0.172062079807758
$ws.Range("J16").Value = 0.188481718551003
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.01597233333333333
$ws.Range("N16").Value = 0.047917
$ws.Range("O16").Value = 0.0034916150995138
$ws.Range("P16").Value = 0.003698320415929645
$ws.Range("Q16").Value = 0.0007304148033333334
$ws.Range("R16").Value = 0.006573733230000001
$ws.Range("S16").Value = 0.0006007745559105164
$ws.Range("T16").Value = 0.0006970657877466798

$ws.Range("A17").Value = 'MuSCs'
$ws.Range("B17").Value = 'Ntrk3'
$ws.Range("C17").Value = 'Ptprf'
$ws.Range("D17").Value = 'ECs'
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.06945950000000001
$ws.Range("H17").Value = 0.138919
$ws.Range("I17").Value = 0.2613458568206203
$ws.Range("J17").Value = 0.1908571459974254
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1255626666666667
$ws.Range("N17").Value = 0.376688
$ws.Range("O17").Value = 0.02744849445093922
$ws.Range("P17").Value = 0.02907345870642374
$ws.Range("Q17").Value = 0.008721520045333336
$ws.Range("R17").Value = 0.05232912027200001
$ws.Range("S17").Value = 0.007173550300716753
$ws.Range("T17").Value = 0.005548877352982033

$ws.Range("A18").Value = 'MuSCs'
$ws.Range("B18").Value = 'Ntrk3'
$ws.Range("C18").Value = 'Ptprf'
$ws.Range("D18").Value = 'FAPs'
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.06945950000000001
$ws.Range("H18").Value = 0.138919
$ws.Range("I18").Value = 0.2613458568206203
$ws.Range("J18").Value = 0.1908571459974254
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 3.393572666666666
$ws.Range("N18").Value = 10.180718
$ws.Range("O18").Value = 0.7418483772500767
$ws.Range("P18").Value = 0.7857661629113346
$ws.Range("Q18").Value = 0.2357158606403333
$ws.Range("R18").Value = 1.414295163842
$ws.Range("S18").Value = 0.1938789997834081
$ws.Range("T18").Value = 0.1499690872746053

$ws.Range("A19").Value = 'MuSCs'
$ws.Range("B19").Value = 'Ntrk3'
$ws.Range("C19").Value = 'Ptprf'
$ws.Range("D19").Value = 'Inflammatory-Mac'
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.06945950000000001
$ws.Range("H19").Value = 0.138919
$ws.Range("I19").Value = 0.2613458568206203
$ws.Range("J19").Value = 0.1908571459974254
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.2723486666666667
$ws.Range("N19").Value = 0.8170460000000001
$ws.Range("O19").Value = 0.05953649332381727
$ws.Range("P19").Value = 0.06306108275880487
$ws.Range("Q19").Value = 0.01891720221233334
$ws.Range("R19").Value = 0.113503213274
$ws.Range("S19").Value = 0.01555961585980817
$ws.Range("T19").Value = 0.01203565827885295

$ws.Range("A20").Value = 'MuSCs'
$ws.Range("B20").Value = 'Ntrk3'
$ws.Range("C20").Value = 'Ptprf'
$ws.Range("D20").Value = 'MuSCs'
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.06945950000000001
$ws.Range("H20").Value = 0.138919
$ws.Range("I20").Value = 0.2613458568206203
$ws.Range("J20").Value = 0.1908571459974254
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 0.7670265000000001
$ws.Range("N20").Value = 1.534053
$ws.Range("O20").Value = 0.167675019875653
$ws.Range("P20").Value = 0.1184009752075072
$ws.Range("Q20").Value = 0.05327727717675001
$ws.Range("R20").Value = 0.213109108707
$ws.Range("S20").Value = 0.04382117173681708
$ws.Range("T20").Value = 0.02259767221141673

$ws.Range("A21").Value = 'MuSCs'
$ws.Range("B21").Value = 'Ntrk3'
$ws.Range("C21").Value = 'Ptprf'
$ws.Range("D21").Value = 'Resolving-Mac'
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.06945950000000001
$ws.Range("H21").Value = 0.138919
$ws.Range("I21").Value = 0.2613458568206203
$ws.Range("J21").Value = 0.1908571459974254
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 0.3333333333333333
$ws.Range("M21").Value = 0.01597233333333333
$ws.Range("N21").Value = 0.047917
$ws.Range("O21").Value = 0.0034916150995138
$ws.Range("P21").Value = 0.003698320415929645
$ws.Range("Q21").Value = 0.001109430287166667
$ws.Range("R21").Value = 0.006656581723000001
$ws.Range("S21").Value = 0.0009125191398702496
$ws.Range("T21").Value = 0.0007058508795683432

$ws.Range("A22").Value = 'Resolving-Mac'
$ws.Range("B22").Value = 'Ntrk3'
$ws.Range("C22").Value = 'Ptprf'
$ws.Range("D22").Value = 'ECs'
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 0.3333333333333333
$ws.Range("G22").Value = 0.1128566666666667
$ws.Range("H22").Value = 0.33857
$ws.Range("I22").Value = 0.4246305004775321
$ws.Range("J22").Value = 0.4651523831898322
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 0.1255626666666667
$ws.Range("N22").Value = 0.376688
$ws.Range("O22").Value = 0.02744849445093922
$ws.Range("P22").Value = 0.02907345870642374
$ws.Range("Q22").Value = 0.01417058401777778
$ws.Range("R22").Value = 0.12753525616
$ws.Range("S22").Value = 0.01165546793605708
$ws.Range("T22").Value = 0.01352358860486418

$ws.Range("A23").Value = 'Resolving-Mac'
$ws.Range("B23").Value = 'Ntrk3'
$ws.Range("C23").Value = 'Ptprf'
$ws.Range("D23").Value = 'FAPs'
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 0.3333333333333333
$ws.Range("G23").Value = 0.1128566666666667
$ws.Range("H23").Value = 0.33857
$ws.Range("I23").Value = 0.4246305004775321
$ws.Range("J23").Value = 0.4651523831898322
$ws.Range("K23").Value = 3
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = 3.393572666666666
$ws.Range("N23").Value = 10.180718
$ws.Range("O23").Value = 0.7418483772500767
$ws.Range("P23").Value = 0.7857661629113346
$ws.Range("Q23").Value = 0.382987299251111
$ws.Range("R23").Value = 3.44688569326
$ws.Range("S23").Value = 0.3150114477101451
$ws.Range("T23").Value = 0.3655010033081372

$ws.Range("A24").Value = 'Resolving-Mac'
$ws.Range("B24").Value = 'Ntrk3'
$ws.Range("C24").Value = 'Ptprf'
$ws.Range("D24").Value = 'Inflammatory-Mac'
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 0.3333333333333333
$ws.Range("G24").Value = 0.1128566666666667
$ws.Range("H24").Value = 0.33857
$ws.Range("I24").Value = 0.4246305004775321
$ws.Range("J24").Value = 0.4651523831898322
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 0.2723486666666667
$ws.Range("N24").Value = 0.8170460000000001
$ws.Range("O24").Value = 0.05953649332381727
$ws.Range("P24").Value = 0.06306108275880487
$ws.Range("Q24").Value = 0.03073636269111111
$ws.Range("R24").Value = 0.27662726422
$ws.Range("S24").Value = 0.02528101095676978
$ws.Range("T24").Value = 0.02933301293178932

$ws.Range("A25").Value = 'Resolving-Mac'
$ws.Range("B25").Value = 'Ntrk3'
$ws.Range("C25").Value = 'Ptprf'
$ws.Range("D25").Value = 'MuSCs'
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 0.3333333333333333
$ws.Range("G25").Value = 0.1128566666666667
$ws.Range("H25").Value = 0.33857
$ws.Range("I25").Value = 0.4246305004775321
$ws.Range("J25").Value = 0.4651523831898322
$ws.Range("K25").Value = 2
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 0.7670265000000001
$ws.Range("N25").Value = 1.534053
$ws.Range("O25").Value = 0.167675019875653
$ws.Range("P25").Value = 0.1184009752075072
$ws.Range("Q25").Value = 0.086564054035
$ws.Range("R25").Value = 0.51938432421
$ws.Range("S25").Value = 0.07119992760737867
$ws.Range("T25").Value = 0.05507449578977219

$ws.Range("A26").Value = 'Resolving-Mac'
$ws.Range("B26").Value = 'Ntrk3'
$ws.Range("C26").Value = 'Ptprf'
$ws.Range("D26").Value = 'Resolving-Mac'
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 0.3333333333333333
$ws.Range("G26").Value = 0.1128566666666667
$ws.Range("H26").Value = 0.33857
$ws.Range("I26").Value = 0.4246305004775321
$ws.Range("J26").Value = 0.4651523831898322
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 0.3333333333333333
$ws.Range("M26").Value = 0.01597233333333333
$ws.Range("N26").Value = 0.047917
$ws.Range("O26").Value = 0.0034916150995138
$ws.Range("P26").Value = 0.003698320415929645
$ws.Range("Q26").Value = 0.001802584298888889
$ws.Range("R26").Value = 0.01622325869
$ws.Range("S26").Value = 0.001482646267181453
$ws.Range("T26").Value = 0.001720282555269286

